$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove all existing hyperlinks up front; they will be re-added fresh below
# (Range.Hyperlinks.Delete() clears every hyperlink on the sheet in this runtime)
$ws.Range("A1").Hyperlinks.Delete()

# Row 2
$ws.Cells.Item(2,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(2,2).Value = '【効率化】Air Tableでデータ収集と工数管理とスムーズにしたい!'
$ws.Cells.Item(2,3).Value = 'システム開発'
$ws.Cells.Item(2,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(2,5).Value = '期限情報なし'
$ws.Cells.Item(2,6).Value = 'https://www.lancers.jp/work/detail/5416639'
$ws.Cells.Item(2,7).Value = 388
$ws.Cells.Item(2,8).Value = '🔥AI,Ai ◆効率化 ◇管理'

# Row 3
$ws.Cells.Item(3,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(3,2).Value = '医療機関向けAIアプリとLINEの連携開発を支援してくださるAIエンジニア募集(AI/バックエンド)'
$ws.Cells.Item(3,3).Value = 'システム開発'
$ws.Cells.Item(3,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(3,5).Value = '期限情報なし'
$ws.Cells.Item(3,6).Value = 'https://www.lancers.jp/work/detail/5416301'
$ws.Cells.Item(3,7).Value = 385
$ws.Cells.Item(3,8).Value = '🔥AI,Ai ◆開発 ◇アプリ'

# Row 4
$ws.Cells.Item(4,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(4,2).Value = '生成AI使用可 【急募】モバイルフィットネストラッカーアプリ開発者募集'
$ws.Cells.Item(4,3).Value = 'システム開発'
$ws.Cells.Item(4,4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(4,5).Value = '期限情報なし'
$ws.Cells.Item(4,6).Value = 'https://www.lancers.jp/work/detail/5417041'
$ws.Cells.Item(4,7).Value = 378
$ws.Cells.Item(4,8).Value = '🔥AI,Ai ◆開発 ◇アプリ'

# Row 5
$ws.Cells.Item(5,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(5,2).Value = '大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5416307'
$ws.Cells.Item(5,7).Value = 378
$ws.Cells.Item(5,8).Value = '🔥AI,Ai ◆効率化'

# Row 6
$ws.Cells.Item(6,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(6,2).Value = 'Kintone × SharePoint × PowerAutomate連携業務システム(AI連携)'
$ws.Cells.Item(6,3).Value = 'システム開発'
$ws.Cells.Item(6,4).Value = '5,000,000 円 ~ / 固定'
$ws.Cells.Item(6,5).Value = '期限情報なし'
$ws.Cells.Item(6,6).Value = 'https://www.lancers.jp/work/detail/5416528'
$ws.Cells.Item(6,7).Value = 325
$ws.Cells.Item(6,8).Value = '🔥AI,Ai'

# Row 7
$ws.Cells.Item(7,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(7,2).Value = '詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発'
$ws.Cells.Item(7,3).Value = 'システム開発'
$ws.Cells.Item(7,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(7,5).Value = '期限情報なし'
$ws.Cells.Item(7,6).Value = 'https://www.lancers.jp/work/detail/5377709'
$ws.Cells.Item(7,7).Value = 245
$ws.Cells.Item(7,8).Value = '🔥Next.js ◆開発,Node.js ◇アプリ'

# Row 8
$ws.Cells.Item(8,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(8,2).Value = '<Next.js、バックエンド開発> ガントチャートアプリの改修製造'
$ws.Cells.Item(8,3).Value = 'システム開発'
$ws.Cells.Item(8,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(8,5).Value = '期限情報なし'
$ws.Cells.Item(8,6).Value = 'https://www.lancers.jp/work/detail/5379158'
$ws.Cells.Item(8,7).Value = 225
$ws.Cells.Item(8,8).Value = '🔥Next.js ◆開発 ◇アプリ'

# Row 9
$ws.Cells.Item(9,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(9,2).Value = '【日本人限定/継続案件】Node.jsエンジニア募集(スクレイピング機能開発)'
$ws.Cells.Item(9,3).Value = 'システム開発'
$ws.Cells.Item(9,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(9,5).Value = '期限情報なし'
$ws.Cells.Item(9,6).Value = 'https://www.lancers.jp/work/detail/5416511'
$ws.Cells.Item(9,7).Value = 155
$ws.Cells.Item(9,8).Value = '◆開発,Node.js'

# Row 10
$ws.Cells.Item(10,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(10,2).Value = 'Flutterなどハイブリッドアプリによる業務アプリの開発(スマートウォッチ)'
$ws.Cells.Item(10,3).Value = 'システム開発'
$ws.Cells.Item(10,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(10,5).Value = '期限情報なし'
$ws.Cells.Item(10,6).Value = 'https://www.lancers.jp/work/detail/5379176'
$ws.Cells.Item(10,7).Value = 100
$ws.Cells.Item(10,8).Value = '◆開発 ◇アプリ'

# Row 11
$ws.Cells.Item(11,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(11,2).Value = '【急募】Laravelでのバックエンド開発:管理画面機能やDB管理・ポイント機能などの開発'
$ws.Cells.Item(11,3).Value = 'システム開発'
$ws.Cells.Item(11,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(11,5).Value = '期限情報なし'
$ws.Cells.Item(11,6).Value = 'https://www.lancers.jp/work/detail/5416675'
$ws.Cells.Item(11,7).Value = 100
$ws.Cells.Item(11,8).Value = '◆開発 ◇管理'

# Row 12
$ws.Cells.Item(12,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(12,2).Value = '完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします'
$ws.Cells.Item(12,3).Value = 'システム開発'
$ws.Cells.Item(12,4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(12,5).Value = '期限情報なし'
$ws.Cells.Item(12,6).Value = 'https://www.lancers.jp/work/detail/5416665'
$ws.Cells.Item(12,7).Value = 70
$ws.Cells.Item(12,8).Value = '◆効率化'

# Row 13
$ws.Cells.Item(13,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(13,2).Value = 'IB報酬を得るための高性能EA開発依頼'
$ws.Cells.Item(13,3).Value = 'システム開発'
$ws.Cells.Item(13,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(13,5).Value = '期限情報なし'
$ws.Cells.Item(13,6).Value = 'https://www.lancers.jp/work/detail/5416508'
$ws.Cells.Item(13,7).Value = 68
$ws.Cells.Item(13,8).Value = '◆開発'

# Row 14
$ws.Cells.Item(14,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(14,2).Value = '【急募】全国店舗をGoogleマップで表示するWPプラグイン開発'
$ws.Cells.Item(14,3).Value = 'システム開発'
$ws.Cells.Item(14,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(14,5).Value = '期限情報なし'
$ws.Cells.Item(14,6).Value = 'https://www.lancers.jp/work/detail/5416539'
$ws.Cells.Item(14,7).Value = 63
$ws.Cells.Item(14,8).Value = '◆開発'

# Row 15
$ws.Cells.Item(15,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(15,2).Value = '開発メンター兼プロジェクト推進パートナー募集!あなたの技術と経験で、開発チームを次のステージへ。'
$ws.Cells.Item(15,3).Value = 'システム開発'
$ws.Cells.Item(15,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(15,5).Value = '期限情報なし'
$ws.Cells.Item(15,6).Value = 'https://www.lancers.jp/work/detail/5416656'
$ws.Cells.Item(15,7).Value = 60
$ws.Cells.Item(15,8).Value = '◆開発'

# Row 16
$ws.Cells.Item(16,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(16,2).Value = '【急募】GASを使った顧客管理スプレッドシートの作成・改修依頼'
$ws.Cells.Item(16,3).Value = 'システム開発'
$ws.Cells.Item(16,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(16,5).Value = '期限情報なし'
$ws.Cells.Item(16,6).Value = 'https://www.lancers.jp/work/detail/5416338'
$ws.Cells.Item(16,7).Value = 33
$ws.Cells.Item(16,8).Value = '◇管理'

# Row 17
$ws.Cells.Item(17,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(17,2).Value = 'サイトスピードが遅く サイトスピードを速くしたい ワードプレス'
$ws.Cells.Item(17,3).Value = 'システム開発'
$ws.Cells.Item(17,4).Value = '20,000 円 ~ 30,000 円 / 募集期間 3 日、取引期間 0 日'
$ws.Cells.Item(17,5).Value = '期限情報なし'
$ws.Cells.Item(17,6).Value = 'https://www.lancers.jp/work/detail/5416402'
$ws.Cells.Item(17,7).Value = 30
$ws.Cells.Item(17,8).Value = '◇サイト'

# Row 18
$ws.Cells.Item(18,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(18,2).Value = 'OR(operations research)にて最適化の仕組みの構築(社内常駐)'
$ws.Cells.Item(18,3).Value = 'システム開発'
$ws.Cells.Item(18,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(18,5).Value = '期限情報なし'
$ws.Cells.Item(18,6).Value = 'https://www.lancers.jp/work/detail/5372984'
$ws.Cells.Item(18,7).Value = 25
$ws.Cells.Item(18,8).ClearContents()

# Row 19
$ws.Cells.Item(19,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(19,2).Value = 'OR(operations research)にて最適化の仕組みの構築 (リモート)'
$ws.Cells.Item(19,3).Value = 'システム開発'
$ws.Cells.Item(19,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(19,5).Value = '期限情報なし'
$ws.Cells.Item(19,6).Value = 'https://www.lancers.jp/work/detail/5367840'
$ws.Cells.Item(19,7).Value = 25
$ws.Cells.Item(19,8).ClearContents()

# Row 20
$ws.Cells.Item(20,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(20,2).Value = '〖リモート可〗Delphiエンジニア募集'
$ws.Cells.Item(20,3).Value = 'システム開発'
$ws.Cells.Item(20,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(20,5).Value = '期限情報なし'
$ws.Cells.Item(20,6).Value = 'https://www.lancers.jp/work/detail/5341051'
$ws.Cells.Item(20,7).Value = 25
$ws.Cells.Item(20,8).ClearContents()

# Row 21
$ws.Cells.Item(21,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(21,2).Value = '初回 【継続あり】Microsoft PL-300/400/600 資格試験向け問題集作成'
$ws.Cells.Item(21,3).Value = 'システム開発'
$ws.Cells.Item(21,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(21,5).Value = '期限情報なし'
$ws.Cells.Item(21,6).Value = 'https://www.lancers.jp/work/detail/5411149'
$ws.Cells.Item(21,7).Value = 18
$ws.Cells.Item(21,8).ClearContents()

# Row 22
$ws.Cells.Item(22,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(22,2).Value = '【急募】キントーン見積書をエクセルに変換してくれる方'
$ws.Cells.Item(22,3).Value = 'システム開発'
$ws.Cells.Item(22,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(22,5).Value = '期限情報なし'
$ws.Cells.Item(22,6).Value = 'https://www.lancers.jp/work/detail/5416819'
$ws.Cells.Item(22,7).Value = 18
$ws.Cells.Item(22,8).ClearContents()

# Row 23
$ws.Cells.Item(23,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(23,2).Value = '【継続案件あり】AWSに精通しているインフラエンジニアを募集します'
$ws.Cells.Item(23,3).Value = 'システム開発'
$ws.Cells.Item(23,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(23,5).Value = '期限情報なし'
$ws.Cells.Item(23,6).Value = 'https://www.lancers.jp/work/detail/5416510'
$ws.Cells.Item(23,7).Value = 10
$ws.Cells.Item(23,8).ClearContents()

# Row 24
$ws.Cells.Item(24,1).Value = '2025-10-20 12:36:29'
$ws.Cells.Item(24,2).Value = '【急募】エクセルマクロの組み方を教えてください!'
$ws.Cells.Item(24,3).Value = 'システム開発'
$ws.Cells.Item(24,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(24,5).Value = '期限情報なし'
$ws.Cells.Item(24,6).Value = 'https://www.lancers.jp/work/detail/5416433'
$ws.Cells.Item(24,7).Value = 10
$ws.Cells.Item(24,8).ClearContents()

# Re-add hyperlinks for every row, in order, matching the URL text already in column F
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), 'https://www.lancers.jp/work/detail/5416639')
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), 'https://www.lancers.jp/work/detail/5416301')
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), 'https://www.lancers.jp/work/detail/5417041')
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), 'https://www.lancers.jp/work/detail/5416307')
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), 'https://www.lancers.jp/work/detail/5416528')
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), 'https://www.lancers.jp/work/detail/5377709')
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), 'https://www.lancers.jp/work/detail/5379158')
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), 'https://www.lancers.jp/work/detail/5416511')
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), 'https://www.lancers.jp/work/detail/5379176')
$ws.Hyperlinks.Add($ws.Cells.Item(11,6), 'https://www.lancers.jp/work/detail/5416675')
$ws.Hyperlinks.Add($ws.Cells.Item(12,6), 'https://www.lancers.jp/work/detail/5416665')
$ws.Hyperlinks.Add($ws.Cells.Item(13,6), 'https://www.lancers.jp/work/detail/5416508')
$ws.Hyperlinks.Add($ws.Cells.Item(14,6), 'https://www.lancers.jp/work/detail/5416539')
$ws.Hyperlinks.Add($ws.Cells.Item(15,6), 'https://www.lancers.jp/work/detail/5416656')
$ws.Hyperlinks.Add($ws.Cells.Item(16,6), 'https://www.lancers.jp/work/detail/5416338')
$ws.Hyperlinks.Add($ws.Cells.Item(17,6), 'https://www.lancers.jp/work/detail/5416402')
$ws.Hyperlinks.Add($ws.Cells.Item(18,6), 'https://www.lancers.jp/work/detail/5372984')
$ws.Hyperlinks.Add($ws.Cells.Item(19,6), 'https://www.lancers.jp/work/detail/5367840')
$ws.Hyperlinks.Add($ws.Cells.Item(20,6), 'https://www.lancers.jp/work/detail/5341051')
$ws.Hyperlinks.Add($ws.Cells.Item(21,6), 'https://www.lancers.jp/work/detail/5411149')
$ws.Hyperlinks.Add($ws.Cells.Item(22,6), 'https://www.lancers.jp/work/detail/5416819')
$ws.Hyperlinks.Add($ws.Cells.Item(23,6), 'https://www.lancers.jp/work/detail/5416510')
$ws.Hyperlinks.Add($ws.Cells.Item(24,6), 'https://www.lancers.jp/work/detail/5416433')

# Column H width: 17 -> 27 characters
$ws.Columns.Item(8).ColumnWidth = 26.166666666666668
